$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.950.44'
$ws.Range("E2").Value = '  +1.12%  '

$ws.Range("D3").Value = '3.931.71'
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '488.15'
$ws.Range("E5").Value = '  +0.64%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.73'
$ws.Range("E6").Value = '  +0.52%  '

$ws.Range("E7").Value = '  -0.74%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.735'
$ws.Range("E9").Value = '  +1.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.176'
$ws.Range("E10").Value = '  +3.32%  '

$ws.Range("E11").Value = '  -4.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.02'
$ws.Range("E12").Value = '  +0.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.45'
$ws.Range("E13").Value = '  -1.32%  '

$ws.Range("D14").Value = '4.558.12'
$ws.Range("E14").Value = '  +0.21%  '

$ws.Range("D15").Value = '3.932.13'
$ws.Range("E15").Value = '  -0.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.22'
$ws.Range("E16").Value = '  -4.26%  '

$ws.Range("E17").Value = '  -0.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.94'
$ws.Range("E18").Value = '  +0.20%  '

$ws.Range("E19").Value = '  +3.06%  '

$ws.Range("D20").Value = '69.009.01'
$ws.Range("E20").Value = '  +1.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.14'
$ws.Range("E21").Value = '  -2.32%  '

$ws.Range("E22").Value = '  +2.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.54'
$ws.Range("E23").Value = '  -2.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.49'
$ws.Range("E24").Value = '  +16.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '89.35'
$ws.Range("E25").Value = '  +0.78%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.74'
$ws.Range("E26").Value = '  +3.65%  '

$ws.Range("E27").Value = '  -2.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.11'
$ws.Range("E28").Value = '  -4.41%  '

$ws.Range("E29").Value = '  -3.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '709.92'
$ws.Range("E30").Value = '  +2.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.52'
$ws.Range("E31").Value = '  +0.62%  '

$ws.Range("E32").Value = '  +1.05%  '

$ws.Range("E33").Value = '  +1.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.481'
$ws.Range("E34").Value = '  +30.43%  '

$ws.Range("D35").Value = '0.0₃0888'
$ws.Range("E35").Value = '  -6.23%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '61.77'
$ws.Range("E36").Value = '  +4.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.05'
$ws.Range("E37").Value = '  +7.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '40.79'
$ws.Range("E38").Value = '  -2.05%  '

$ws.Range("E39").Value = '  -0.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0490'
$ws.Range("E42").Value = '  +2.24%  '

$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.96'
$ws.Range("E43").Value = '  +2.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.07'
$ws.Range("E44").Value = '  -1.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.00'
$ws.Range("E45").Value = '  +1.13%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.36'
$ws.Range("E47").Value = '  +5.69%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0364'
$ws.Range("E48").Value = '  +9.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.01'
$ws.Range("E49").Value = '  +5.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.39'
$ws.Range("E50").Value = '  -1.01%  '

$ws.Range("E51").Value = '  -3.09%  '
